# Update "想去人数" (want-to-go count) figures in the 南宁-漫展信息 workbook.
# Sheet "展览" (Exhibitions), sheet "演出" (Performances) and sheet "全部类型"
# (All types, a union of the two) each carry their own copy of these figures
# in column F, keyed by row.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 612
$wsExhibit.Range("F4").Value = 528
$wsExhibit.Range("F5").Value = 516
$wsExhibit.Range("F6").Value = 290
$wsExhibit.Range("F7").Value = 2614
$wsExhibit.Range("F8").Value = 446
$wsExhibit.Range("F9").Value = 7170
$wsExhibit.Range("F10").Value = 190
$wsExhibit.Range("F11").Value = 448
$wsExhibit.Range("F12").Value = 16
$wsExhibit.Range("F13").Value = 146

# --- Sheet: 演出 (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 2

# --- Sheet: 全部类型 (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 612
$wsAll.Range("F4").Value = 528
$wsAll.Range("F5").Value = 516
$wsAll.Range("F6").Value = 290
$wsAll.Range("F9").Value = 2614
$wsAll.Range("F10").Value = 446
$wsAll.Range("F11").Value = 7170
$wsAll.Range("F12").Value = 190
$wsAll.Range("F13").Value = 448
$wsAll.Range("F14").Value = 16
$wsAll.Range("F15").Value = 2
$wsAll.Range("F17").Value = 146
